# Update the "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Mapping of row -> (old value, new value) is driven by the published diff.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1152
    3  = 585
    5  = 33
    6  = 142
    10 = 5213
    11 = 4789
    15 = 48
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
